$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Periodo Mora" (period) labels in column E for rows 16-22 so
# that the six most-recent periods (2407-2412) precede 2501, and update the
# "Valor Mora" amounts in column F to match the row each period now occupies.
$ws.Range("E16").Value = "2407"
$ws.Range("F16").Value = 62580

$ws.Range("E17").Value = "2408"
$ws.Range("F17").Value = 62580

$ws.Range("E18").Value = "2409"
$ws.Range("F18").Value = 62580

$ws.Range("E19").Value = "2410"
$ws.Range("F19").Value = 62580

$ws.Range("E20").Value = "2411"
$ws.Range("F20").Value = 62580

$ws.Range("E21").Value = "2412"
$ws.Range("F21").Value = 62580

$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 60494
